# Auto-generated edit script to update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells keep their original text formatting
# (values look numeric, so force Text format before assigning)
$textCells = @('D5', 'D6', 'D11', 'D14', 'D18', 'D21', 'D22', 'D24', 'D26', 'D27', 'D28', 'D31', 'D34', 'D35', 'D37', 'D40', 'D42', 'D43', 'D44', 'D47', 'D48', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '70.879.00'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '3.805.01'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('D5').Value = '697.61'
$ws.Range('E5').Value = '  +11.08%  '
$ws.Range('D6').Value = '173.50'
$ws.Range('E6').Value = '  +4.55%  '
$ws.Range('D7').Value = '3.802.76'
$ws.Range('E7').Value = '  +0.89%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('D11').Value = '7.58'
$ws.Range('E11').Value = '  +11.82%  '
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('E13').Value = '  +3.08%  '
$ws.Range('D14').Value = '36.26'
$ws.Range('E14').Value = '  +3.95%  '
$ws.Range('D15').Value = '4.450.97'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').Value = '3.806.38'
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').Value = '70.921.07'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').Value = '17.77'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').Value = '11.39'
$ws.Range('E21').Value = '  +19.71%  '
$ws.Range('D22').Value = '479.21'
$ws.Range('E22').Value = '  +3.47%  '
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('D24').Value = '83.65'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  +1.23%  '
$ws.Range('D26').Value = '12.37'
$ws.Range('E26').Value = '  +3.28%  '
$ws.Range('D27').Value = '2.15'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').Value = '10.33'
$ws.Range('E28').Value = '  +2.72%  '
$ws.Range('D29').Value = '3.957.66'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = '3.11'
$ws.Range('E31').Value = '  +16.13%  '
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('E33').Value = '  +5.57%  '
$ws.Range('D34').Value = '29.59'
$ws.Range('E34').Value = '  +3.75%  '
$ws.Range('D35').Value = '0.178'
$ws.Range('E35').Value = '  +4.83%  '
$ws.Range('E36').Value = '  +2.27%  '
$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').Value = '0.997'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('B38').Value = 'RenzoRestakedETH'
$ws.Range('C38').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D38').Value = '3.756.58'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('D40').Value = '3.54'
$ws.Range('E40').Value = '  +7.20%  '
$ws.Range('E41').Value = '  +3.04%  '
$ws.Range('D42').Value = '0.000335'
$ws.Range('E42').Value = '  +25.95%  '
$ws.Range('D43').Value = '2.20'
$ws.Range('E43').Value = '  +13.40%  '
$ws.Range('D44').Value = '0.966'
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = '49.48'
$ws.Range('E47').Value = '  +6.12%  '
$ws.Range('D48').Value = '160.26'
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').Value = '45.01'
$ws.Range('E50').Value = '  +4.41%  '
$ws.Range('E51').Value = '  +1.42%  '
